$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '54.191.99'
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").Value = '2.263.14'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").Value = '495.72'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").Value = '128.66'
$ws.Range("E6").Value = '  +0.84%  '
$ws.Range("E7").Value = '  +0.40%  '
$ws.Range("E8").Value = '  -0.55%  '
$ws.Range("D9").Value = '0.0951'
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("E10").Value = '  +0.87%  '
$ws.Range("E11").Value = '  +2.83%  '
$ws.Range("E12").Value = '  +4.36%  '
$ws.Range("D13").Value = '22.89'
$ws.Range("E13").Value = '  +5.30%  '
$ws.Range("D14").Value = '2.663.62'
$ws.Range("E14").Value = '  +0.26%  '
$ws.Range("D15").Value = '54.170.39'
$ws.Range("E15").Value = '  +0.35%  '
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").Value = '2.261.42'
$ws.Range("E17").Value = '  -1.42%  '
$ws.Range("E18").Value = '  +2.13%  '
$ws.Range("E19").Value = '  +0.68%  '
$ws.Range("D20").Value = '302.65'
$ws.Range("E20").Value = '  +0.66%  '
$ws.Range("D21").Value = '6.32'
$ws.Range("E21").Value = '  -1.52%  '
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.44%  '
$ws.Range("D23").Value = '60.60'
$ws.Range("E23").Value = '  -2.65%  '
$ws.Range("D24").Value = '0.999'
$ws.Range("E24").Value = '  -2.02%  '
$ws.Range("E25").Value = '  +0.10%  '
$ws.Range("E26").Value = '  +2.96%  '
$ws.Range("D27").Value = '171.78'
$ws.Range("E27").Value = '  +2.13%  '
$ws.Range("E28").Value = '  -0.37%  '
$ws.Range("E29").Value = '  +1.64%  '
$ws.Range("D30").Value = '0.0₃0688'
$ws.Range("E30").Value = '  +0.29%  '
$ws.Range("D31").Value = '1.07'
$ws.Range("E31").Value = '  +0.52%  '
$ws.Range("E32").Value = '  +0.13%  '
$ws.Range("D33").Value = '17.76'
$ws.Range("E33").Value = '  +0.26%  '
$ws.Range("E34").Value = '  +0.12%  '
$ws.Range("D35").Value = '0.941'
$ws.Range("E35").Value = '  +3.33%  '
$ws.Range("D36").Value = '1.19'
$ws.Range("E36").Value = '  +0.44%  '
$ws.Range("D37").Value = '3.69'
$ws.Range("E37").Value = '  +0.24%  '
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("E39").Value = '  -0.33%  '
$ws.Range("E40").Value = '  +0.09%  '
$ws.Range("D41").Value = '4.80'
$ws.Range("E41").Value = '  -3.12%  '
$ws.Range("D42").Value = '124.44'
$ws.Range("E42").Value = '  -1.65%  '
$ws.Range("E43").Value = '  +1.55%  '
$ws.Range("D44").Value = '0.0893'
$ws.Range("E44").Value = '  +0.85%  '
$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("D46").Value = '240.52'
$ws.Range("E46").Value = '  +1.27%  '
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("D48").Value = '0.0203'
$ws.Range("E48").Value = '  +0.56%  '
$ws.Range("E49").Value = '  +0.87%  '
$ws.Range("D50").Value = '16.08'
$ws.Range("E50").Value = '  -1.05%  '
$ws.Range("E51").Value = '  -0.27%  '
